# Table 13 "Catch dependency factors per country" — update CATCH.DEP values
# for Cod / Hake per country (no TAC, no RECOVERY, no sp dependency).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $new) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    # Exclude the trailing cell-mark character so only the visible text
    # of the cell is replaced; this avoids Find/Replace matching an
    # unrelated occurrence of the same old value elsewhere in the doc.
    $inner = $d.Range($r.Start, $r.End - 1)
    $inner.Text = $new
}

# Row 2 = BE
Set-CellText $t 2 2 "0.970"
Set-CellText $t 2 3 "0.996"

# Row 3 = DK
Set-CellText $t 3 2 "0.959"
Set-CellText $t 3 3 "0.995"

# Row 4 = DE
Set-CellText $t 4 2 "0.764"
Set-CellText $t 4 3 "0.999"

# Row 5 = EE
Set-CellText $t 5 2 "0.948"

# Row 6 = IE
Set-CellText $t 6 2 "0.989"

# Row 7 = ES
Set-CellText $t 7 2 "0.871"
Set-CellText $t 7 3 "0.916"

# Row 8 = FR
Set-CellText $t 8 2 "0.977"
Set-CellText $t 8 3 "0.838"

# Row 9 = LV
Set-CellText $t 9 2 "0.714"

# Row 10 = LT
Set-CellText $t 10 2 "0.464"

# Row 11 = NL
Set-CellText $t 11 2 "1.000"

# Row 12 = PL
Set-CellText $t 12 2 "0.401"

# Row 13 = PT
Set-CellText $t 13 2 "0.839"
Set-CellText $t 13 3 "0.939"

# Row 14 = FI
Set-CellText $t 14 2 "0.921"

# Row 15 = SE
Set-CellText $t 15 2 "0.722"
Set-CellText $t 15 3 "0.996"
